# Apply updated market-price data to the Lamia_Profits workbook.
# Values sourced from a scheduled market-board refresh; only the
# price/profit columns (H-N) on a handful of rows per job sheet change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 46723.91
$ws.Range("I28").Value = 67474.07000000001
$ws.Range("J28").Value = 2259.2856
$ws.Range("K28").Value = 67474.07000000001
$ws.Range("L28").Value = 2259.2856
$ws.Range("M28").Value = -66989.07000000001
$ws.Range("N28").Value = -3229.2856
$ws.Range("H33").Value = 379.45456
$ws.Range("I33").Value = 375.22223
$ws.Range("K33").Value = 375.22223
$ws.Range("M33").Value = -146.22223
$ws.Range("H34").Value = 4040.7144
$ws.Range("J34").Value = 7998.3335
$ws.Range("L34").Value = 7998.3335
$ws.Range("N34").Value = -8404.333500000001
$ws.Range("H36").Value = 4040.7144
$ws.Range("J36").Value = 7998.3335
$ws.Range("L36").Value = 7998.3335
$ws.Range("N36").Value = -9428.333500000001
$ws.Range("H40").Value = 4096.8237
$ws.Range("I40").Value = 3633.25
$ws.Range("J40").Value = 4349.6816
$ws.Range("K40").Value = 3633.25
$ws.Range("L40").Value = 4349.6816
$ws.Range("M40").Value = -3458.25
$ws.Range("N40").Value = -4699.6816
$ws.Range("H53").Value = 2113.1667
$ws.Range("I53").Value = 400
$ws.Range("K53").Value = 400
$ws.Range("M53").Value = 237
$ws.Range("H62").Value = 4704.3335
$ws.Range("I62").Value = 3919.875
$ws.Range("J62").Value = 6273.25
$ws.Range("K62").Value = 3919.875
$ws.Range("L62").Value = 6273.25
$ws.Range("M62").Value = -3295.875
$ws.Range("N62").Value = -7521.25
$ws.Range("H65").Value = 4704.3335
$ws.Range("I65").Value = 3919.875
$ws.Range("J65").Value = 6273.25
$ws.Range("K65").Value = 19599.375
$ws.Range("L65").Value = 31366.25
$ws.Range("M65").Value = -16479.375
$ws.Range("N65").Value = -37606.25
$ws.Range("H106").Value = 7386.0303
$ws.Range("I106").Value = 2510.8096
$ws.Range("K106").Value = 2510.8096
$ws.Range("M106").Value = -1879.8096
$ws.Range("H113").Value = 7532.8887
$ws.Range("J113").Value = 8685.143
$ws.Range("L113").Value = 8685.143
$ws.Range("N113").Value = -15193.143
$ws.Range("H132").Value = 1308.9762
$ws.Range("I132").Value = 1107.0256
$ws.Range("K132").Value = 3321.0768
$ws.Range("M132").Value = -791.0767999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 302.4
$ws.Range("I5").Value = 254.5
$ws.Range("K5").Value = 254.5
$ws.Range("M5").Value = -142.5
$ws.Range("H74").Value = 10418252
$ws.Range("I74").Value = 12347279
$ws.Range("K74").Value = 12347279
$ws.Range("M74").Value = -12346405
$ws.Range("H77").Value = 10418252
$ws.Range("I77").Value = 12347279
$ws.Range("K77").Value = 61736395
$ws.Range("M77").Value = -61732027
$ws.Range("H131").Value = 89288.5
$ws.Range("J131").Value = 89288.5
$ws.Range("L131").Value = 89288.5
$ws.Range("N131").Value = -99368.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 302.4
$ws.Range("I4").Value = 254.5
$ws.Range("K4").Value = 254.5
$ws.Range("M4").Value = -139.5
$ws.Range("H20").Value = 3886.842
$ws.Range("I20").Value = 2738.6365
$ws.Range("K20").Value = 2738.6365
$ws.Range("M20").Value = -2491.6365
$ws.Range("H22").Value = 433.27274
$ws.Range("I22").Value = 426.7
$ws.Range("K22").Value = 426.7
$ws.Range("M22").Value = -253.7
$ws.Range("H105").Value = 13707.167
$ws.Range("I105").Value = 13710.883
$ws.Range("K105").Value = 13710.883
$ws.Range("M105").Value = -11963.883

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37254
$ws.Range("I31").Value = 4385.75
$ws.Range("K31").Value = 4385.75
$ws.Range("M31").Value = -4090.75
$ws.Range("H34").Value = 37254
$ws.Range("I34").Value = 4385.75
$ws.Range("K34").Value = 4385.75
$ws.Range("M34").Value = -4183.75
$ws.Range("H132").Value = 2016.6086
$ws.Range("I132").Value = 1218.4
$ws.Range("K132").Value = 3655.2
$ws.Range("M132").Value = -1125.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 461.8889
$ws.Range("J11").Value = 42
$ws.Range("L11").Value = 126
$ws.Range("N11").Value = -406
$ws.Range("H23").Value = 446.8
$ws.Range("I23").Value = 200
$ws.Range("J23").Value = 508.5
$ws.Range("K23").Value = 600
$ws.Range("L23").Value = 1525.5
$ws.Range("M23").Value = -365
$ws.Range("N23").Value = -1995.5
$ws.Range("H38").Value = 49.625
$ws.Range("I38").Value = 75.75
$ws.Range("J38").Value = 23.5
$ws.Range("K38").Value = 227.25
$ws.Range("L38").Value = 70.5
$ws.Range("M38").Value = 119.75
$ws.Range("N38").Value = -764.5
$ws.Range("H88").Value = 8064.5
$ws.Range("J88").Value = 7073.7144
$ws.Range("L88").Value = 21221.1432
$ws.Range("N88").Value = -22077.1432
$ws.Range("H91").Value = 8064.5
$ws.Range("J91").Value = 7073.7144
$ws.Range("L91").Value = 21221.1432
$ws.Range("N91").Value = -24185.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 5099.75
$ws.Range("I113").Value = 4625.5835
$ws.Range("J113").Value = 6522.25
$ws.Range("K113").Value = 4625.5835
$ws.Range("L113").Value = 6522.25
$ws.Range("M113").Value = -2455.5835
$ws.Range("N113").Value = -10862.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7626.65
$ws.Range("I40").Value = 6235.2666
$ws.Range("K40").Value = 6235.2666
$ws.Range("M40").Value = -6099.2666
$ws.Range("H68").Value = 3723.36
$ws.Range("I68").Value = 3090.0908
$ws.Range("J68").Value = 8367.333000000001
$ws.Range("K68").Value = 3090.0908
$ws.Range("L68").Value = 8367.333000000001
$ws.Range("M68").Value = -2341.0908
$ws.Range("N68").Value = -9865.333000000001
$ws.Range("H71").Value = 3723.36
$ws.Range("I71").Value = 3090.0908
$ws.Range("J71").Value = 8367.333000000001
$ws.Range("K71").Value = 15450.454
$ws.Range("L71").Value = 41836.665
$ws.Range("M71").Value = -11706.454
$ws.Range("N71").Value = -49324.665
$ws.Range("H82").Value = 3315.5833
$ws.Range("I82").Value = 3203.25
$ws.Range("J82").Value = 3540.25
$ws.Range("K82").Value = 3203.25
$ws.Range("L82").Value = 3540.25
$ws.Range("M82").Value = -2842.25
$ws.Range("N82").Value = -4262.25
$ws.Range("H85").Value = 3315.5833
$ws.Range("I85").Value = 3203.25
$ws.Range("J85").Value = 3540.25
$ws.Range("K85").Value = 3203.25
$ws.Range("L85").Value = 3540.25
$ws.Range("M85").Value = -1955.25
$ws.Range("N85").Value = -6036.25
$ws.Range("H133").Value = 60326
$ws.Range("J133").Value = 60326
$ws.Range("L133").Value = 60326
$ws.Range("N133").Value = -65386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H81").Value = 6791.4
$ws.Range("I81").Value = 3152.5
$ws.Range("J81").Value = 12249.75
$ws.Range("K81").Value = 6305
$ws.Range("L81").Value = 24499.5
$ws.Range("M81").Value = -5244
$ws.Range("N81").Value = -26621.5
$ws.Range("H84").Value = 6791.4
$ws.Range("I84").Value = 3152.5
$ws.Range("J84").Value = 12249.75
$ws.Range("K84").Value = 31525
$ws.Range("L84").Value = 122497.5
$ws.Range("M84").Value = -26221
$ws.Range("N84").Value = -133105.5
$ws.Range("H101").Value = 35000
$ws.Range("J101").Value = 35000
$ws.Range("L101").Value = 35000
$ws.Range("N101").Value = -41490
$ws.Range("H102").Value = 74993
$ws.Range("J102").Value = 74993
$ws.Range("L102").Value = 74993
$ws.Range("N102").Value = -81483
$ws.Range("H103").Value = 32499.5
$ws.Range("J103").Value = 32499.5
$ws.Range("L103").Value = 32499.5
$ws.Range("N103").Value = -34843.5
$ws.Range("H105").Value = 29997
$ws.Range("J105").Value = 29997
$ws.Range("L105").Value = 29997
$ws.Range("N105").Value = -36985
$ws.Range("H132").Value = 6947.3335
$ws.Range("I132").Value = 6259.073
$ws.Range("K132").Value = 18777.219
$ws.Range("M132").Value = -16247.219

Write-Output "Applied all Lamia_Profits market-data updates"
